$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.771.64'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.06%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.491.59'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.73'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.19'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.93%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.76%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.214'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.72%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.660'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.46%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.42'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.33%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000306'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.85%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.59'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.062.70'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.66%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '618.59'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.75%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.914.28'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.66'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.63%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.83'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.63%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.502.65'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.31%  '

# Row 20
$ws.Range("E20").Value = '  +0.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.989'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.28%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.85'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.44%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.48'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +11.58%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.63'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.05'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.89%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.04'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.20%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.15%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.88'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +6.12%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.27'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +6.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.13'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.01%  '

# Row 31
$ws.Range("E31").Value = '  +12.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.59'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.37%  '

# Row 33
$ws.Range("E33").Value = '  -0.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.18'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.52%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.726.89'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.17%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '524.08'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.65%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.07'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.40%  '

# Row 38
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0794'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.82%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.390'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.81%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.58'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.97%  '

# Row 42
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.63'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.74%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.16%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0462'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.43%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.75%  '

# Row 46
$ws.Range("E46").Value = '  +1.90%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.76%  '

# Row 48
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.43%  '

# Row 49
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.10%  '

# Row 50
$ws.Range("E50").Value = '  -1.90%  '

# Row 51
$ws.Range("B51").Value = 'OceanProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.34'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.51%  '
